# Applies two edits described by the diff:
#  1) Insert a new, empty ListParagraph-styled paragraph (ind left=294,
#     rPr lang=en-US) immediately before the "[4] Feedback
#     (Suggestion/Query):" heading paragraph.
#  2) Remove the <w:lastRenderedPageBreak/> marker from the run that
#     contains "Company admin can give feedback to super admin any time."

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: insert the new empty paragraph before the "[4] Feedback
# (Suggestion/Query):" paragraph.
# ---------------------------------------------------------------------
$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "[4] Feedback (Suggestion/Query):`r") {
        $headingPara = $p
        break
    }
}

if ($headingPara -ne $null) {
    $headingStart = $headingPara.Range.Start
    # Zero-length range right at the start of the heading paragraph -
    # InsertXML there inserts a new paragraph before it instead of
    # clobbering the heading paragraph's own content.
    $insertionPoint = $d.Range($headingStart, $headingStart)

    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:ind w:left="294"/>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '</w:pPr>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $null = $insertionPoint.InsertXML($newParaXml)
}

# ---------------------------------------------------------------------
# Edit 2: drop <w:lastRenderedPageBreak/> from the "Company admin can
# give feedback to super admin any time." run.
# ---------------------------------------------------------------------
$feedbackPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Company admin can give feedback to super admin any time.`r") {
        $feedbackPara = $p
        break
    }
}

if ($feedbackPara -ne $null) {
    $paraRange = $feedbackPara.Range
    # Exclude the trailing paragraph mark so only the run is replaced;
    # this leaves the owning <w:p>'s own attributes untouched.
    $runRange = $d.Range($paraRange.Start, $paraRange.End - 1)

    $runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
        '<w:t>Company admin can give feedback to super admin any time.</w:t>' +
        '</w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $null = $runRange.InsertXML($runXml)
}
